# Export template update: the report title (merged cell A1:F1) is switched
# from the hard-coded Vietnamese caption "DANH SACH GIANG VIEN CUA KHOA
# CONG NGHE THONG TIN" to the "{{Name}}" template placeholder, matching the
# other Mustache-style tags ({{Items.Id}}, {{Items.Name}}, ...) used across
# this export sheet. The stale D15 selection left over from editing is also
# moved back onto the title row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "{{Name}}"
$ws.Range("A1:F1").Select()
